# Auto-generated script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.016.06"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -2.17%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.502.97"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -3.33%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.03%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'551.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -3.42%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'147.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -4.66%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.02%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.617"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.95%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.499.43"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -3.40%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.108"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -8.81%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -1.49%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'5.39"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -7.68%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.357"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -5.89%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'26.25"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -6.80%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'2.951.44"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -3.44%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'61.868.21"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -2.15%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.0000164"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -7.95%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.501.12"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -3.27%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'11.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -6.73%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'7.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -5.81%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -7.37%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'322.47"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -5.56%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D24").Value = "'63.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -4.86%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'1.75"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -4.13%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -5.48%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'2.623.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -3.27%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("B28").Value = "'Fetch.AI"
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").Value = "'1.50"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -3.45%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("B29").Value = "'Binance-PegBSC-USD"
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = "'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.26%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'537.29"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -6.66%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'8.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -7.71%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'7.69"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -2.21%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -6.11%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.90"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -7.53%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.58"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -8.14%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'5.92"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -9.12%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'4.90"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -9.76%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -0.08%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.379"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -5.62%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'18.56"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -5.78%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'144.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -6.64%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.00%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'1.70"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -8.08%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'40.39"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -2.05%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'2.31"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -6.29%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'149.55"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -4.10%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'3.59"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -8.03%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'20.88"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -9.48%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.0537"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -8.42%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0954"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -4.85%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -5.62%  "
$ws.Range("E51").Style = "Normal"

Write-Output "Applied 95 cell updates"
